$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: num_customers 24 -> 25, retention_rate recalculated (C22/D22)
$ws.Range("C22").Value = 25
$ws.Range("E22").Value = 0.009419743782969104

# Row 27: num_customers 36 -> 37, retention_rate recalculated (C27/D27)
$ws.Range("C27").Value = 37
$ws.Range("E27").Value = 0.01642984014209592

# Row 34: num_customers 62 -> 68, retention_rate recalculated (C34/D34)
$ws.Range("C34").Value = 68
$ws.Range("E34").Value = 0.0301418439716312

# Row 36: num_customers 99 -> 102, retention_rate recalculated (C36/D36)
$ws.Range("C36").Value = 102
$ws.Range("E36").Value = 0.05284974093264249

# Row 37: num_customers and cohort_size 604 -> 630 (retention_rate stays 1)
$ws.Range("C37").Value = 630
$ws.Range("D37").Value = 630
